# Completed heat transfer coefficient calculation
# Add the new "Heat Transfer Coefficient" column (R) to the Calculations sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in R1, matching the style used by the other header cells (Q1).
$ws.Range("R1").Value = "Heat Transfer Coefficient"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# New calculated data for each experiment row.
$ws.Range("R2").Value = 93.367340989409044
$ws.Range("R3").Value = 112.42574902742849
$ws.Range("R4").Value = 82.643192133033367
$ws.Range("R5").Value = 65.742888212295099
